$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "27+42="
$t.Cell(1,2).Range.Text = "74-36="
$t.Cell(1,3).Range.Text = "73-66="
$t.Cell(1,4).Range.Text = "49-48="
$t.Cell(1,5).Range.Text = "13+11="
$t.Cell(2,1).Range.Text = "2+88="
$t.Cell(2,2).Range.Text = "28+28="
$t.Cell(2,3).Range.Text = "30+44="
$t.Cell(2,4).Range.Text = "52-24="
$t.Cell(2,5).Range.Text = "59-14="
$t.Cell(3,1).Range.Text = "74-72="
$t.Cell(3,2).Range.Text = "39+12="
$t.Cell(3,3).Range.Text = "65+16="
$t.Cell(3,4).Range.Text = "94-71="
$t.Cell(3,5).Range.Text = "91-67="
$t.Cell(4,1).Range.Text = "78-41="
$t.Cell(4,2).Range.Text = "12+78="
$t.Cell(4,3).Range.Text = "54+36="
$t.Cell(4,4).Range.Text = "48-24="
$t.Cell(4,5).Range.Text = "71-41="
$t.Cell(5,1).Range.Text = "31+17="
$t.Cell(5,2).Range.Text = "86+1="
$t.Cell(5,3).Range.Text = "30+32="
$t.Cell(5,4).Range.Text = "26+70="
$t.Cell(5,5).Range.Text = "75+19="
$t.Cell(6,1).Range.Text = "72-8="
$t.Cell(6,2).Range.Text = "58-52="
$t.Cell(6,3).Range.Text = "87-11="
$t.Cell(6,4).Range.Text = "55-11="
$t.Cell(6,5).Range.Text = "65+34="
$t.Cell(7,1).Range.Text = "10+47="
$t.Cell(7,2).Range.Text = "81-31="
$t.Cell(7,3).Range.Text = "69+10="
$t.Cell(7,4).Range.Text = "41+23="
$t.Cell(7,5).Range.Text = "42+48="
$t.Cell(8,1).Range.Text = "42+13="
$t.Cell(8,2).Range.Text = "23+75="
$t.Cell(8,3).Range.Text = "23+67="
$t.Cell(8,4).Range.Text = "41+45="
$t.Cell(8,5).Range.Text = "96-60="
$t.Cell(9,1).Range.Text = "79+2="
$t.Cell(9,2).Range.Text = "19+62="
$t.Cell(9,3).Range.Text = "38+42="
$t.Cell(9,4).Range.Text = "22+40="
$t.Cell(9,5).Range.Text = "94-15="
$t.Cell(10,1).Range.Text = "37+10="
$t.Cell(10,2).Range.Text = "3+29="
$t.Cell(10,3).Range.Text = "89-22="
$t.Cell(10,4).Range.Text = "56+23="
$t.Cell(10,5).Range.Text = "33+43="
$t.Cell(11,1).Range.Text = "56+13="
$t.Cell(11,2).Range.Text = "44+36="
$t.Cell(11,3).Range.Text = "72-52="
$t.Cell(11,4).Range.Text = "29+42="
$t.Cell(11,5).Range.Text = "98-87="
$t.Cell(12,1).Range.Text = "83-41="
$t.Cell(12,2).Range.Text = "46+53="
$t.Cell(12,3).Range.Text = "69-15="
$t.Cell(12,4).Range.Text = "4+5="
$t.Cell(12,5).Range.Text = "42-40="
$t.Cell(13,1).Range.Text = "86-76="
$t.Cell(13,2).Range.Text = "77+21="
$t.Cell(13,3).Range.Text = "34+18="
$t.Cell(13,4).Range.Text = "65-55="
$t.Cell(13,5).Range.Text = "54-54="
$t.Cell(14,1).Range.Text = "92-89="
$t.Cell(14,2).Range.Text = "8+1="
$t.Cell(14,3).Range.Text = "67-16="
$t.Cell(14,4).Range.Text = "6+11="
$t.Cell(14,5).Range.Text = "30-24="
$t.Cell(15,1).Range.Text = "20+56="
$t.Cell(15,2).Range.Text = "39+15="
$t.Cell(15,3).Range.Text = "33+32="
$t.Cell(15,4).Range.Text = "59+19="
$t.Cell(15,5).Range.Text = "74-71="
$t.Cell(16,1).Range.Text = "96-17="
$t.Cell(16,2).Range.Text = "15-6="
$t.Cell(16,3).Range.Text = "89-56="
$t.Cell(16,4).Range.Text = "2+92="
$t.Cell(16,5).Range.Text = "2+81="
$t.Cell(17,1).Range.Text = "39-26="
$t.Cell(17,2).Range.Text = "34+25="
$t.Cell(17,3).Range.Text = "98-32="
$t.Cell(17,4).Range.Text = "18+26="
$t.Cell(17,5).Range.Text = "33+63="
$t.Cell(18,1).Range.Text = "54+30="
$t.Cell(18,2).Range.Text = "66-6="
$t.Cell(18,3).Range.Text = "19+43="
$t.Cell(18,4).Range.Text = "10-7="
$t.Cell(18,5).Range.Text = "95-74="
$t.Cell(19,1).Range.Text = "10+46="
$t.Cell(19,2).Range.Text = "51-46="
$t.Cell(19,3).Range.Text = "20+38="
$t.Cell(19,4).Range.Text = "19+65="
$t.Cell(19,5).Range.Text = "92-23="
$t.Cell(20,1).Range.Text = "7+15="
$t.Cell(20,2).Range.Text = "89-42="
$t.Cell(20,3).Range.Text = "6+5="
$t.Cell(20,4).Range.Text = "8+56="
$t.Cell(20,5).Range.Text = "18+7="
